$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet
$ws.Range("Z9").Interior.TintAndShade = 0.39997558519241921
$ws.Range("Z9").Interior.ThemeColor = 10
